$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# About sheet: add a region label + "last updated" date in the header row
# ---------------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")
$about.Range("B1").Value = "Colorado"
$about.Range("C1").Value = 45271
$about.Range("C1").NumberFormat = "mm-dd-yy"

# ---------------------------------------------------------------------------
# FPIEBP sheet: flip the production/imports/exports balancing priority order
# from (production=1, imports=3, exports=2) to (exports=1, imports=2,
# production=3) for every fuel row except hydrogen (row 22), which keeps its
# original order.
# ---------------------------------------------------------------------------
$fpiebp = $wb.Worksheets.Item("FPIEBP")

$rows = 3,4,5,9,10,11,12,13,14,17,18,19,20,21
foreach ($r in $rows) {
    $fpiebp.Range("B$r").Value = 3
    $fpiebp.Range("C$r").Value = 2
    $fpiebp.Range("D$r").Value = 1
}

# E9 was a leftover empty styled cell with no value; drop it entirely so the
# sheet's used range shrinks back from E22 to D22.
$fpiebp.Range("E9").Clear()

# The workbook was left with the FPIEBP tab active and the selection parked
# just past the data table.
$fpiebp.Activate() | Out-Null
$fpiebp.Range("F4").Select() | Out-Null
